$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct small mapping errors in the "land" column (C) and a few others,
# aligning them with the correct regioinvent country code in column A.

$ws.Range("E2").Value = "RER"

$ws.Range("C3").Value = "AE"

$ws.Range("C41").Value = "CV"
$ws.Range("C42").Value = "CY"
$ws.Range("C50").Value = "EG"

$ws.Range("B66").Value = "HK"

$ws.Range("C72").Value = "IL"
$ws.Range("C75").Value = "IS"
$ws.Range("C78").Value = "JO"

$ws.Range("B83").Value = "KI"

$ws.Range("C86").Value = "KW"
$ws.Range("C87").Value = "KY"
$ws.Range("C90").Value = "LB"
$ws.Range("C93").Value = "LS"
$ws.Range("C97").Value = "LY"
$ws.Range("C109").Value = "MT"
$ws.Range("C124").Value = "OM"

$ws.Range("D132").Value = "PS"
$ws.Range("E132").Value = "PS"

$ws.Range("C135").Value = "QA"
$ws.Range("C141").Value = "SA"
$ws.Range("C153").Value = "TJ"
$ws.Range("C154").Value = "TL"
$ws.Range("C163").Value = "UZ"
$ws.Range("C167").Value = "YE"

# Update the selected cell to match the author's last active selection.
$ws.Range("E19").Select()
